# Reproduce the commit: add a new worksheet "SO4_excel" to the workbook,
# placed after the existing "soil" sheet, containing a full copy of the
# "SO4" sheet (same data + formatting), and make the new sheet the
# active / selected tab.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("SO4")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy "SO4" to a new sheet placed right after the last existing sheet
# (i.e. after "soil"), mirroring Excel's own Worksheet.Copy semantics.
$source.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "SO4_excel"

# Make the freshly added sheet the active tab, like in the authored file.
$newSheet.Activate()
